$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Refresh the access tokens for consumer22 (Sheet1) and provider23 (Sheet2)
$ws1.Range("F2").Value = "eyJ0eXAiOiJKV1QiLCJhbGciOiJIUzI1NiJ9.eyJ0b2tlbl90eXBlIjoiYWNjZXNzIiwiZXhwIjoxNjUxNTA0MjkxLCJqdGkiOiIwNjBiM2QzODdlNDA0NDZkOTc3ZGRmYTA4OGUxNjQyYSIsInVzZXJfaWQiOjc1Miwicm9sZSI6IkMiLCJ1c2VybmFtZSI6ImNvbnN1bWVyMjIiLCJlbWFpbCI6ImNvbnN1bWVyMjJAZ21haWwuY29tIn0.8T5RUFtgVVxf-gGlU6RT8zt5G7GlPdBS_8r_CyGmTeA"
$ws2.Range("F2").Value = "eyJ0eXAiOiJKV1QiLCJhbGciOiJIUzI1NiJ9.eyJ0b2tlbl90eXBlIjoiYWNjZXNzIiwiZXhwIjoxNjUxNTA0MzA3LCJqdGkiOiIyZGQ1ZjBhYzcwYzY0ODBjYTIyNmNmNjMyZjdlYjY0YyIsInVzZXJfaWQiOjc1MCwicm9sZSI6IlAiLCJ1c2VybmFtZSI6InByb3ZpZGVyMjMiLCJlbWFpbCI6InByb3ZpZGVyMjNAZ21haWwuY29tIn0.biui-AfgGNMZvvy_hF5wL4hFzI0xv3i0eyTaVjxF5Lo"

# Add a new consumer vehicle record below the existing data on Sheet1
$ws1.Range("A5").Value = "TS08UF4343"
$ws1.Range("B5").Value = "'2"
$ws1.Range("C5").Value = "Black SUV"
$ws1.Range("D5").Value = "Mahindra"

$ws1.Range("A6").Value = "'37"

# Match the author's final selection
$ws1.Activate()
$ws1.Range("A6").Select()
